$wb = $excel.ActiveWorkbook

# New row (row 99) data for each of the 4 worksheets, in workbook tab order:
# ROW35-FE-LIFTER, ROW35-MID-LIFTER, ROW02-FE-LIFTER, ROW02-MID-LIFTER

$rowsData = @(
    @{ A = 45773.95323103009;  B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x56"; E = "0xd"; F = 400; G = [double]"5.68631262647114e+23"; H = 342; I = 13 },
    @{ A = 45773.81087542824;  B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x56"; E = "0xe"; F = 400; G = [double]"5.68631262647114e+23"; H = 342; I = 14 },
    @{ A = 45773.95240554398;  B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x56"; E = "0x3"; F = 400; G = [double]"5.68631262647114e+23"; H = 342; I = 3 },
    @{ A = 45774.01796109954;  B = "0x01,0x90"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x01,0x56"; E = "0x3"; F = 400; G = [double]"9.85046333984776e+23"; H = 342; I = 3 }
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$i - 1]
    $newRow = 99

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}

Write-Output "done"
